# Update computed bus voltage magnitudes (vm_pu) for the 380 kV case.
# Source sheet: Code/Results/Cases/Case_2_0/res_bus/vm_pu.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.067708051918255
$ws.Range("D2").Value = 1.068928237496619
$ws.Range("E2").Value = 1.069053404763189
$ws.Range("F2").Value = 1.079712618211368
$ws.Range("I2").Value = 1.046923539255988
$ws.Range("J2").Value = 1.072651575665976
$ws.Range("K2").Value = 1.071632108104883
$ws.Range("L2").Value = 1.071756940574727
$ws.Range("M2").Value = 1.082387950457629
$ws.Range("N2").Value = 1.074174863938278

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.06937090949603
$ws.Range("D3").Value = 1.070218652064131
$ws.Range("E3").Value = 1.070353580390201
$ws.Range("F3").Value = 1.081095009730812
$ws.Range("I3").Value = 1.047323147505314
$ws.Range("J3").Value = 1.073967555187215
$ws.Range("K3").Value = 1.072737417855915
$ws.Range("L3").Value = 1.072872011591538
$ws.Range("M3").Value = 1.083587092872194
$ws.Range("N3").Value = 1.075492712301381

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.070445069764911
$ws.Range("D4").Value = 1.071051882981439
$ws.Range("E4").Value = 1.071194027146775
$ws.Range("F4").Value = 1.081988253966468
$ws.Range("I4").Value = 1.047579570399759
$ws.Range("J4").Value = 1.074816893526217
$ws.Range("K4").Value = 1.073450321434756
$ws.Range("L4").Value = 1.073592131003729
$ws.Range("M4").Value = 1.08436122436976
$ws.Range("N4").Value = 1.076343256798245

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.070896221237328
$ws.Range("D5").Value = 1.071401760237615
$ws.Range("E5").Value = 1.071547152266114
$ws.Range("F5").Value = 1.082363479998296
$ws.Range("I5").Value = 1.047686858412523
$ws.Range("J5").Value = 1.075173439213574
$ws.Range("K5").Value = 1.073749480455778
$ws.Range("L5").Value = 1.073894537922482
$ws.Range("M5").Value = 1.084686246062578
$ws.Range("N5").Value = 1.076700308821374

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.070971946896782
$ws.Range("D6").Value = 1.071460482145191
$ws.Range("E6").Value = 1.071606432040918
$ws.Range("F6").Value = 1.082426465063571
$ws.Range("I6").Value = 1.04770484258397
$ws.Range("J6").Value = 1.075233274760548
$ws.Range("K6").Value = 1.073799678758811
$ws.Range("L6").Value = 1.073945294102752
$ws.Range("M6").Value = 1.084740794
$ws.Range("N6").Value = 1.076760229341688

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.070451099730793
$ws.Range("D7").Value = 1.07105655967181
$ws.Range("E7").Value = 1.071198746394566
$ws.Range("F7").Value = 1.081993268900206
$ws.Range("I7").Value = 1.047581005996201
$ws.Range("J7").Value = 1.074821659724417
$ws.Range("K7").Value = 1.073454320947326
$ws.Range("L7").Value = 1.073596173073943
$ws.Range("M7").Value = 1.084365568979566
$ws.Range("N7").Value = 1.076348029764993

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.068270405348188
$ws.Range("D8").Value = 1.06936470633728
$ws.Range("E8").Value = 1.069492985438633
$ws.Range("F8").Value = 1.080180067923372
$ws.Range("I8").Value = 1.047059035570376
$ws.Range("J8").Value = 1.073096774915998
$ws.Range("K8").Value = 1.072006133967143
$ws.Range("L8").Value = 1.072134078412958
$ws.Range("M8").Value = 1.082793581454154
$ws.Range("N8").Value = 1.074620695422299

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.064413349446932
$ws.Range("D9").Value = 1.066369708339966
$ws.Range("E9").Value = 1.06648042047575
$ws.Range("F9").Value = 1.076975057524141
$ws.Range("I9").Value = 1.046122669240503
$ws.Range("J9").Value = 1.070040197810593
$ws.Range("K9").Value = 1.069436298147461
$ws.Range("L9").Value = 1.069546671278922
$ws.Range("M9").Value = 1.080009518278086
$ws.Range("N9").Value = 1.07155977762681

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.061831660844746
$ws.Range("D10").Value = 1.064363377269776
$ws.Range("E10").Value = 1.06446711241895
$ws.Range("F10").Value = 1.074831289047093
$ws.Range("I10").Value = 1.045487105181285
$ws.Range("J10").Value = 1.067990485629749
$ws.Range("K10").Value = 1.067710606350466
$ws.Range("L10").Value = 1.067813992300684
$ws.Range("M10").Value = 1.078143653254999
$ws.Range("N10").Value = 1.069507154619564

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.060711170754222
$ws.Range("D11").Value = 1.063492224659312
$ws.Range("E11").Value = 1.063594076797129
$ws.Range("F11").Value = 1.073901235197097
$ws.Range("I11").Value = 1.045209176368649
$ws.Range("J11").Value = 1.067099982929723
$ws.Range("K11").Value = 1.066960315243
$ws.Range("L11").Value = 1.067061811049876
$ws.Range("M11").Value = 1.077333294272696
$ws.Range("N11").Value = 1.068615387303562

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.060294567532558
$ws.Range("D12").Value = 1.063168271107772
$ws.Range("E12").Value = 1.063269596108124
$ws.Range("F12").Value = 1.073555494838488
$ws.Range("I12").Value = 1.045105528466022
$ws.Range("J12").Value = 1.066768755746279
$ws.Range("K12").Value = 1.066681156842747
$ws.Range("L12").Value = 1.066782122370307
$ws.Range("M12").Value = 1.077031918462874
$ws.Range("N12").Value = 1.068283689739525

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.060383948788869
$ws.Range("D13").Value = 1.063237777063815
$ws.Range("E13").Value = 1.063339207318917
$ws.Range("F13").Value = 1.073629669986055
$ws.Range("I13").Value = 1.045127780025003
$ws.Range("J13").Value = 1.06683982583785
$ws.Range("K13").Value = 1.066741058521492
$ws.Range("L13").Value = 1.066842130001404
$ws.Range("M13").Value = 1.077096581574399
$ws.Range("N13").Value = 1.068354860758778

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06067674250156
$ws.Range("D14").Value = 1.063465454127641
$ws.Range("E14").Value = 1.063567259169659
$ws.Range("F14").Value = 1.073872661889302
$ws.Range("I14").Value = 1.045200617240256
$ws.Range("J14").Value = 1.067072612921393
$ws.Range("K14").Value = 1.066937249492472
$ws.Range("L14").Value = 1.067038697954386
$ws.Range("M14").Value = 1.077308390096257
$ws.Range("N14").Value = 1.068587978426681

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.060857088555573
$ws.Range("D15").Value = 1.063605684373885
$ws.Range("E15").Value = 1.063707743205027
$ws.Range("F15").Value = 1.074022340271897
$ws.Range("I15").Value = 1.045245439857369
$ws.Range("J15").Value = 1.067215980169736
$ws.Range("K15").Value = 1.067058067122488
$ws.Range("L15").Value = 1.067159770658156
$ws.Range("M15").Value = 1.077438842736958
$ws.Range("N15").Value = 1.068731549272961

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.061905968206337
$ws.Range("D16").Value = 1.064421141559985
$ws.Range("E16").Value = 1.064525025768187
$ws.Range("F16").Value = 1.074892975278785
$ws.Range("I16").Value = 1.045505492698846
$ws.Range("J16").Value = 1.06804952210376
$ws.Range("K16").Value = 1.067760335661057
$ws.Range("L16").Value = 1.067863871079192
$ws.Range("M16").Value = 1.078197382314863
$ws.Range("N16").Value = 1.06956627493214

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.062563197373918
$ws.Range("D17").Value = 1.064932008378525
$ws.Range("E17").Value = 1.065037343215465
$ws.Range("F17").Value = 1.07543861707451
$ws.Range("I17").Value = 1.04566788498726
$ws.Range("J17").Value = 1.068571580725203
$ws.Range("K17").Value = 1.068200026670522
$ws.Range("L17").Value = 1.068305016011242
$ws.Range("M17").Value = 1.078672538890298
$ws.Range("N17").Value = 1.070089074936706

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.06294629782716
$ws.Range("D18").Value = 1.065229757549164
$ws.Range("E18").Value = 1.065336048231246
$ws.Range("F18").Value = 1.075756708757923
$ws.Range("I18").Value = 1.045762342927293
$ws.Range("J18").Value = 1.068875803372909
$ws.Range("K18").Value = 1.068456196753397
$ws.Range("L18").Value = 1.06856214360404
$ws.Range("M18").Value = 1.078949456028411
$ws.Range("N18").Value = 1.070393729615467

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.063076883116727
$ws.Range("D19").Value = 1.065331243426897
$ws.Range("E19").Value = 1.065437878597908
$ws.Range("F19").Value = 1.075865140887532
$ws.Range("I19").Value = 1.045794506179052
$ws.Range("J19").Value = 1.068979487327367
$ws.Range("K19").Value = 1.068543494459592
$ws.Range("L19").Value = 1.068649786325853
$ws.Range("M19").Value = 1.079043838218017
$ws.Range("N19").Value = 1.070497560813033

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.062492708873413
$ws.Range("D20").Value = 1.064877221137454
$ws.Range("E20").Value = 1.064982388977536
$ws.Range("F20").Value = 1.075380092719079
$ws.Range("I20").Value = 1.045650489031391
$ws.Range("J20").Value = 1.068515598356315
$ws.Range("K20").Value = 1.06815288247842
$ws.Range("L20").Value = 1.068257704496019
$ws.Range("M20").Value = 1.078621583310619
$ws.Range("N20").Value = 1.070033013066433

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.060590533284429
$ws.Range("D21").Value = 1.06339841915244
$ws.Range("E21").Value = 1.063500109062201
$ws.Range("F21").Value = 1.073801114542591
$ws.Range("I21").Value = 1.045179179916898
$ws.Range("J21").Value = 1.067004075560061
$ws.Range("K21").Value = 1.066879489093086
$ws.Range("L21").Value = 1.06698082179063
$ws.Range("M21").Value = 1.077246028108529
$ws.Range("N21").Value = 1.068519343734435

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.059392220331494
$ws.Range("D22").Value = 1.062466500338649
$ws.Range("E22").Value = 1.06256700105157
$ws.Range("F22").Value = 1.072806741777669
$ws.Range("I22").Value = 1.044880459367315
$ws.Range("J22").Value = 1.066051084880585
$ws.Range("K22").Value = 1.066076150978548
$ws.Range("L22").Value = 1.06617628442474
$ws.Range("M22").Value = 1.076379003458344
$ws.Range("N22").Value = 1.067564999698884

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.060027694643052
$ws.Range("D23").Value = 1.062960733553526
$ws.Range("E23").Value = 1.063061769687161
$ws.Range("F23").Value = 1.073334032673022
$ws.Range("I23").Value = 1.045039044420877
$ws.Range("J23").Value = 1.066556536373087
$ws.Range("K23").Value = 1.06650227477573
$ws.Range("L23").Value = 1.066602949275465
$ws.Range("M23").Value = 1.076838836837868
$ws.Range("N23").Value = 1.068071168990482

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.062524560365842
$ws.Range("D24").Value = 1.064901977846742
$ws.Range("E24").Value = 1.065007220806561
$ws.Range("F24").Value = 1.075406537887657
$ws.Range("I24").Value = 1.045658350326951
$ws.Range("J24").Value = 1.068540895260806
$ws.Range("K24").Value = 1.068174185818862
$ws.Range("L24").Value = 1.068279083104409
$ws.Range("M24").Value = 1.078644608665868
$ws.Range("N24").Value = 1.07005834589543

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.065412262212148
$ws.Range("D25").Value = 1.067145658996216
$ws.Range("E25").Value = 1.067260083194876
$ws.Range("F25").Value = 1.077804848593734
$ws.Range("I25").Value = 1.046366725231529
$ws.Range("J25").Value = 1.070832474094726
$ws.Range("K25").Value = 1.070102830143406
$ws.Range("L25").Value = 1.07021691741875
$ws.Range("M25").Value = 1.080730967426512
$ws.Range("N25").Value = 1.072353179034142
